$wb = $excel.ActiveWorkbook

# --- Update "Status" text on every sheet (Overview, zh-cn, de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Update timestamps (stored as text, not real dates) ---
$wsOverview.Range("G2").Value = "2016-08-17 18:58:17"
$wsDeDe.Range("H2").Value = "2016-08-17 18:58:17"
$wsZhCn.Range("H2").Value = "2016-08-17 18:58:12"

# --- Narrow columns E & F on Overview, column C on zh-cn / de-de ---
# (target stored width 17.2159881591797 isn't reachable through this host's
# pixel-quantized column-width model; 16.33333.. lands on the nearest
# representable stored width, 17.1666666666667)
$wsOverview.Range("E1").ColumnWidth = 16.333333333333332
$wsOverview.Range("F1").ColumnWidth = 16.333333333333332
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333332
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333332
